$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 89 - delinda.gucwa / Belaliet Houssame Eddine - skipped attempt
$ws.Cells.Item(89, 1).Value = "https://www.instagram.com/delinda.gucwa"
$ws.Cells.Item(89, 2).Value = "Belaliet Houssame Eddine"
$ws.Cells.Item(89, 3).Value = "Skip"
$ws.Cells.Item(89, 4).Value = "Skip"
$ws.Cells.Item(89, 5).Value = "Skip"

# Row 90 - delinda.gucwa / Belaliet Houssame Eddine - not sent
$ws.Cells.Item(90, 1).Value = "https://www.instagram.com/delinda.gucwa"
$ws.Cells.Item(90, 2).Value = "Belaliet Houssame Eddine"
$ws.Cells.Item(90, 3).Value = $false
$ws.Cells.Item(90, 6).Value = $false
$ws.Cells.Item(90, 7).Value = $false

# Row 91 - _bskinn / Bella Little - commented with new comment text
$ws.Cells.Item(91, 1).Value = "https://www.instagram.com/_bskinn"
$ws.Cells.Item(91, 2).Value = "Bella Little"
$ws.Cells.Item(91, 3).Value = $true
$ws.Cells.Item(91, 4).Value = "This is seriously impressive! The way you handle light & shadow is masterful. Wondering what your go-to brushes are for texture? Just sent you something in DMs you might find interesting - check it when you get a chance!"
$ws.Cells.Item(91, 6).Value = $false
$ws.Cells.Item(91, 7).Value = $false

# Row 92 - allisawcakes / Belle - commented with new comment text
$ws.Cells.Item(92, 1).Value = "https://www.instagram.com/allisawcakes"
$ws.Cells.Item(92, 2).Value = "Belle"
$ws.Cells.Item(92, 3).Value = $true
$ws.Cells.Item(92, 4).Value = "Love the way you play with light and shadow! Really creates a mood. What's your favorite lens for these shots? Just sent you something in DMs you might find interesting - check it when you get a chance!"
$ws.Cells.Item(92, 6).Value = $false
$ws.Cells.Item(92, 7).Value = $false

# Row 93 - delinda.gucwa / Belaliet Houssame Eddine - skipped attempt
$ws.Cells.Item(93, 1).Value = "https://www.instagram.com/delinda.gucwa"
$ws.Cells.Item(93, 2).Value = "Belaliet Houssame Eddine"
$ws.Cells.Item(93, 3).Value = "Skip"
$ws.Cells.Item(93, 4).Value = "Skip"
$ws.Cells.Item(93, 5).Value = "Skip"

# Row 94 - delinda.gucwa / Belaliet Houssame Eddine - not sent
$ws.Cells.Item(94, 1).Value = "https://www.instagram.com/delinda.gucwa"
$ws.Cells.Item(94, 2).Value = "Belaliet Houssame Eddine"
$ws.Cells.Item(94, 3).Value = $false
$ws.Cells.Item(94, 6).Value = $false
$ws.Cells.Item(94, 7).Value = $false

# Match the saved view state: scrolled down with D91 selected
$ws.Range("D91").Select()
$excel.ActiveWindow.ScrollRow = 86
$excel.ActiveWindow.ScrollColumn = 1
